# Updated cryptos list on Thu Mar 14 03:52:23 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns; rows 34-36 and 51 also
# swap which coin occupies that rank. D-column values that are plain
# decimals are written with a leading apostrophe so Excel keeps them as
# text (matching the original inlineStr cells) instead of coercing them
# to numbers; the quote-prefix style that introduces is reset back to
# "Normal" immediately afterwards so cell formatting is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.216.89'
$ws.Range('E2').Value = '  +1.61%  '
$ws.Range('D3').Value = '3.995.68'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'612.44"
$ws.Range('E5').Value = '  +13.92%  '
$ws.Range('D6').Value = "'167.28"
$ws.Range('E6').Value = '  +11.93%  '
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('D9').Value = "'0.759"
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = "'0.171"
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('D11').Value = "'57.28"
$ws.Range('E11').Value = '  +7.38%  '
$ws.Range('D12').Value = "'0.0000331"
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = "'11.16"
$ws.Range('E13').Value = '  +2.13%  '
$ws.Range('D14').Value = '4.629.79'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '3.992.76'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = "'1.25"
$ws.Range('E16').Value = '  +3.93%  '
$ws.Range('D17').Value = "'14.20"
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').Value = "'20.61"
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '73.067.96'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = "'439.56"
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('E22').Value = '  +16.04%  '
$ws.Range('D23').Value = "'95.97"
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = "'3.39"
$ws.Range('E24').Value = '  -3.50%  '
$ws.Range('D25').Value = "'14.25"
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').Value = "'4.08"
$ws.Range('D27').Value = "'11.13"
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').Value = "'10.55"
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('D29').Value = "'5.96"
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = "'36.12"
$ws.Range('D31').Value = "'7.75"
$ws.Range('E31').Value = '  -8.47%  '
$ws.Range('D32').Value = "'13.75"
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = "'72.31"
$ws.Range('E34').Value = '  +8.76%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = "'48.42"
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('B36').Value = 'PEPE'
$ws.Range('C36').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D36').Value = "'0.0000102"
$ws.Range('E36').Value = '  +16.50%  '
$ws.Range('D37').Value = "'636.08"
$ws.Range('E37').Value = '  -6.09%  '
$ws.Range('D38').Value = "'0.435"
$ws.Range('E38').Value = '  -5.21%  '
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = "'0.999"
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('D42').Value = "'11.11"
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D44').Value = "'3.27"
$ws.Range('E44').Value = '  -5.21%  '
$ws.Range('D45').Value = "'0.0486"
$ws.Range('E45').Value = '  -1.68%  '
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = "'3.45"
$ws.Range('E47').Value = '  +4.05%  '
$ws.Range('D48').Value = "'2.63"
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('D49').Value = "'2.86"
$ws.Range('E49').Value = '  +29.53%  '
$ws.Range('D50').Value = '2.871.30'
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = "'3.04"
$ws.Range('E51').Value = '  -2.72%  '

$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
